$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.689.53"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "'1.743.56"
$ws.Range("E3").Value = "  -5.50%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'237.81"
$ws.Range("E5").Value = "  -8.84%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4934"
$ws.Range("E7").Value = "  -6.42%  "
$ws.Range("D8").Value = "'41.59"
$ws.Range("E8").Value = "  -7.74%  "
$ws.Range("D9").Value = "'0.2420"
$ws.Range("E9").Value = "  -23.40%  "
$ws.Range("D10").Value = "'0.05981"
$ws.Range("E10").Value = "  -12.10%  "
$ws.Range("D11").Value = "'1.742.63"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("D12").Value = "'0.06845"
$ws.Range("E12").Value = "  -12.01%  "
$ws.Range("D13").Value = "'14.78"
$ws.Range("E13").Value = "  -22.82%  "
$ws.Range("D14").Value = "'4.467"
$ws.Range("E14").Value = "  -11.02%  "
$ws.Range("D15").Value = "'77.27"
$ws.Range("E15").Value = "  -12.59%  "
$ws.Range("D16").Value = "'0.5828"
$ws.Range("E16").Value = "  -25.83%  "
$ws.Range("D17").Value = "'0.9988"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'25.730.68"
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("D20").Value = "'11.48"
$ws.Range("E20").Value = "  -17.68%  "
$ws.Range("D21").Value = "'0.000006445"
$ws.Range("E21").Value = "  -18.78%  "
$ws.Range("D22").Value = "'1.961.75"
$ws.Range("E22").Value = "  -5.93%  "
$ws.Range("D23").Value = "'3.965"
$ws.Range("E23").Value = "  -14.00%  "
$ws.Range("D24").Value = "'4.999"
$ws.Range("E24").Value = "  -16.65%  "
$ws.Range("D25").Value = "'7.841"
$ws.Range("E25").Value = "  -16.19%  "
$ws.Range("D26").Value = "'135.96"
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("D27").Value = "'1.473"
$ws.Range("E27").Value = "  -12.48%  "
$ws.Range("D28").Value = "'1.839"
$ws.Range("E28").Value = "  -17.55%  "
$ws.Range("D29").Value = "'14.54"
$ws.Range("E29").Value = "  -14.74%  "
$ws.Range("D30").Value = "'100.76"
$ws.Range("E30").Value = "  -9.30%  "
$ws.Range("D31").Value = "'3.787"
$ws.Range("E31").Value = "  -10.24%  "
$ws.Range("D32").Value = "'0.08113"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("D33").Value = "'3.376"
$ws.Range("E33").Value = "  -17.54%  "
$ws.Range("D34").Value = "'0.04381"
$ws.Range("E34").Value = "  -10.34%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'2.643"
$ws.Range("E36").Value = "  -7.73%  "
$ws.Range("D37").Value = "'1.023"
$ws.Range("E37").Value = "  -10.40%  "
$ws.Range("D38").Value = "'0.6067"
$ws.Range("E38").Value = "  -17.04%  "
$ws.Range("D39").Value = "'2.708"
$ws.Range("E39").Value = "  -12.95%  "
$ws.Range("D40").Value = "'2.073"
$ws.Range("E40").Value = "  -11.07%  "
$ws.Range("D41").Value = "'0.9999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'102.89"
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("D43").Value = "'0.01490"
$ws.Range("E43").Value = "  -14.27%  "
$ws.Range("D44").Value = "'0.7741"
$ws.Range("E44").Value = "  -15.01%  "
$ws.Range("D45").Value = "'5.131"
$ws.Range("E45").Value = "  -13.54%  "
$ws.Range("D46").Value = "'0.3771"
$ws.Range("E46").Value = "  -22.13%  "
$ws.Range("D47").Value = "'0.05110"
$ws.Range("E47").Value = "  -12.38%  "
$ws.Range("D48").Value = "'6.006"
$ws.Range("E48").Value = "  -22.22%  "
$ws.Range("D49").Value = "'0.1071"
$ws.Range("E49").Value = "  -14.11%  "
$ws.Range("D50").Value = "'30.06"
$ws.Range("E50").Value = "  -13.83%  "
$ws.Range("D51").Value = "'52.77"
$ws.Range("E51").Value = "  -12.18%  "
